$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Resultados conocidos para partidos ya jugados (fila 135 y 136) ---
$ws.Range("G135").Value = "Fallo"
$ws.Range("H135").Value = -1

$ws.Range("G136").Value = "Fallo"
$ws.Range("H136").Value = -1

# --- Nuevas filas del tracker (fila 150 y 151) ---
$ws.Range("A150").Value = 14310265
$ws.Range("B150").Value = "'2025-08-09"
$ws.Range("B150").Style = "Normal"
$ws.Range("C150").Value = "Michael Zheng"
$ws.Range("D150").Value = "Garrett Johns"
$ws.Range("E150").Value = "Gana Garrett Johns"
$ws.Range("F150").Value = 3.75

$ws.Range("A151").Value = 14310272
$ws.Range("B151").Value = "'2025-08-09"
$ws.Range("B151").Style = "Normal"
$ws.Range("C151").Value = "August Holmgren"
$ws.Range("D151").Value = "Yu Hsiou Hsu"
$ws.Range("E151").Value = "Gana Yu Hsiou Hsu"
$ws.Range("F151").Value = 2.1
